$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Two new trade rows appended below the header row (the exporter had
# been silently dropping rows 2-3 -- this restores them).
# ---------------------------------------------------------------------

# --- Row 2 : T_0001 ---
$ws.Range("A2").Value = "T_0001"

$ws.Range("B2").Value = 46028.06429717592
$ws.Range("B2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C2").Value = 46028.06483881945
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("D2").Value = "eth-updown-15m-1767663900"
$ws.Range("E2").Value = "Ethereum Up or Down - January 5, 8:45PM-9:00PM ET"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "82305470999346288835007647903068954073105635720820601714459832234291951995392"

$ws.Range("G2").Value = "Up"
$ws.Range("H2").Value = "BUY"
$ws.Range("I2").Value = 0.48
$ws.Range("J2").Value = 0.485
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 0.02500000000000002
$ws.Range("M2").Value = 1.041666666666668
$ws.Range("N2").Value = "MEAN_REVERSION (@ `$0.4850)"
$ws.Range("O2").Value = 46.798516
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 85
$ws.Range("T2").Style = "Normal"
$ws.Range("U2").Style = "Normal"
$ws.Range("V2").Value = "Priority 2"

# --- Row 3 : T_0002 ---
$ws.Range("A3").Value = "T_0002"

$ws.Range("B3").Value = 46028.06469672605
$ws.Range("B3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").Value = 46028.06484028603
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("D3").Value = "eth-updown-15m-1767663900"
$ws.Range("E3").Value = "Ethereum Up or Down - January 5, 8:45PM-9:00PM ET"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "82305470999346288835007647903068954073105635720820601714459832234291951995392"

$ws.Range("G3").Value = "Up"
$ws.Range("H3").Value = "BUY"
$ws.Range("I3").Value = 0.48
$ws.Range("J3").Value = 0.485
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 0.02500000000000002
$ws.Range("M3").Value = 1.041666666666668
$ws.Range("N3").Value = "MEAN_REVERSION (@ `$0.4850)"
$ws.Range("O3").Value = 12.403582
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 85
$ws.Range("T3").Style = "Normal"
$ws.Range("U3").Style = "Normal"
$ws.Range("V3").Value = "Priority 2"
